$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.128.88'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '1.611.89'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.07%  '
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('E8').Value = '  +1.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0620'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.40'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0798'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').Value = '1.836.14'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').Value = '1.610.09'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '26.151.13'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '198.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('E22').Value = '  +2.11%  '
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.60'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.74'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('E28').Value = '  +2.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0475'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.97%  '
$ws.Range('E32').Value = '  +2.33%  '
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('E34').Value = '  +3.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.84%  '
$ws.Range('D36').Value = '1.108.59'
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('E37').Value = '  +1.62%  '
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.34'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.506'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.791'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.798'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.02%  '
$ws.Range('D43').Value = '1.749.05'
$ws.Range('E43').Value = '  +0.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('D46').Value = '0.0₆0109'
$ws.Range('E46').Value = '  +8.86%  '
$ws.Range('E47').Value = '  +9.12%  '
$ws.Range('E48').Value = '  +1.74%  '
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.408'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('E51').Value = '  -0.18%  '
